$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.082572221755981
$ws.Range("B1").Value = 1.864029288291931
$ws.Range("C1").Value = -1
$ws.Range("D1").Value = 2.031523704528809
$ws.Range("E1").Value = 1.145656108856201
